$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (Price, Volume(1h)) new values
$updates = @{
    2  = @("312.55", "8.33%")
    3  = @("32.04", "9.66%")
    4  = @("5.336", "4.61%")
    5  = @("0.07558", "13.12%")
    6  = @("7.832", "6.83%")
    7  = @("3.683", "8.14%")
    8  = @("1.572", "15.29%")
    9  = @("0.9116", "-0.75%")
    10 = @("0.01702", "2,537.84%")
    11 = @("0.1694", "6.77%")
    12 = @("0.07738", "13.46%")
    13 = @("0.08117", "5.91%")
    14 = @("0.03030", "3.20%")
    15 = @("0.09869", "9.84%")
    16 = @("0.001526", "-2.91%")
    17 = @("0.04554", "1.27%")
    18 = @("0.006401", "2.38%")
    19 = @("3.502", "1.46%")
    20 = @("2.241", "0.95%")
    21 = @("0.3269", "1.72%")
    22 = @("0.1317", "0.62%")
    23 = @("4.173", "2.42%")
    25 = @("0.001215", "1.92%")
    26 = @("0.004491", "9.14%")
    27 = @("0.0001297", "8.19%")
    28 = @("0.0001737", "7.42%")
    40 = @("0.04572", "8.19%")
    41 = @("0.007164", "6.43%")
    42 = @("0.1363", "10.10%")
    43 = @("0.002254", "7.87%")
    44 = @("0.01392", "4.03%")
    45 = @("0.00006172", "8.12%")
    46 = @("1.892", "-4.13%")
    47 = @("0.01297", "-0.70%")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dCell = $ws.Range("D$row")
    $dCell.NumberFormat = "@"
    $dCell.Value = $vals[0]

    $eCell = $ws.Range("E$row")
    $eCell.NumberFormat = "@"
    $eCell.Value = $vals[1]
}
